$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ajout des equivalents ports GPIO dans les colonnes "GPIO" (G, J, M)
$updates = @(
    @("G3", 0),
    @("J3", 4),
    @("M3", 7),
    @("G4", 1),
    @("J4", 5),
    @("M4", 8),
    @("G5", 2),
    @("J5", 6),
    @("M5", 9),
    @("G6", 3),
    @("M6", 10),
    @("M7", 11),
    @("M8", 12)
)

foreach ($pair in $updates) {
    $addr = $pair[0]
    $val = $pair[1]
    $rng = $ws.Range($addr)
    $rng.Value = $val
    $rng.Interior.Color = 5296274
}

$ws.Range("H12").Select()
